# Update LR-pairs sheet with the refreshed TPM-based NATMI values.
# Rows 2-7 get new sending/target cluster pairings and recomputed metrics;
# the former rows 8-10 (MuSCs-sending pairs) are removed entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam1"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.861094666666666
$ws.Range("H2").Value = 29.583284
$ws.Range("I2").Value = 0.243709096397741
$ws.Range("J2").Value = 0.2437090963977409
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005673666666666667
$ws.Range("N2").Value = 0.017021
$ws.Range("O2").Value = 0.1234828534325781
$ws.Range("P2").Value = 0.1234828534325781
$ws.Range("Q2").Value = 0.05594856410711111
$ws.Range("R2").Value = 0.5035370769640001
$ws.Range("S2").Value = 0.0300938946306683
$ws.Range("T2").Value = 0.0300938946306683
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam1"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.861094666666666
$ws.Range("H3").Value = 29.583284
$ws.Range("I3").Value = 0.243709096397741
$ws.Range("J3").Value = 0.2437090963977409
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04027333333333333
$ws.Range("N3").Value = 0.12082
$ws.Range("O3").Value = 0.8765171465674219
$ws.Range("P3").Value = 0.876517146567422
$ws.Range("Q3").Value = 0.3971391525422222
$ws.Range("R3").Value = 3.57425237288
$ws.Range("S3").Value = 0.2136152017670727
$ws.Range("T3").Value = 0.2136152017670727
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icam1"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 29.35342966666667
$ws.Range("H4").Value = 88.060289
$ws.Range("I4").Value = 0.7254466225154019
$ws.Range("J4").Value = 0.7254466225154018
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005673666666666667
$ws.Range("N4").Value = 0.017021
$ws.Range("O4").Value = 0.1234828534325781
$ws.Range("P4").Value = 0.1234828534325781
$ws.Range("Q4").Value = 0.1665415754521111
$ws.Range("R4").Value = 1.498874179069
$ws.Range("S4").Value = 0.0895802189612282
$ws.Range("T4").Value = 0.0895802189612282
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam1"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.35342966666667
$ws.Range("H5").Value = 88.060289
$ws.Range("I5").Value = 0.7254466225154019
$ws.Range("J5").Value = 0.7254466225154018
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04027333333333333
$ws.Range("N5").Value = 0.12082
$ws.Range("O5").Value = 0.8765171465674219
$ws.Range("P5").Value = 0.876517146567422
$ws.Range("Q5").Value = 1.182160457442222
$ws.Range("R5").Value = 10.63944411698
$ws.Range("S5").Value = 0.6358664035541737
$ws.Range("T5").Value = 0.6358664035541737
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Icam1"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.248038666666667
$ws.Range("H6").Value = 3.744116
$ws.Range("I6").Value = 0.03084428108685718
$ws.Range("J6").Value = 0.03084428108685716
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005673666666666667
$ws.Range("N6").Value = 0.017021
$ws.Range("O6").Value = 0.1234828534325781
$ws.Range("P6").Value = 0.1234828534325781
$ws.Range("Q6").Value = 0.007080955381777779
$ws.Range("R6").Value = 0.06372859843600001
$ws.Range("S6").Value = 0.003808739840681626
$ws.Range("T6").Value = 0.003808739840681625
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Icam1"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.248038666666667
$ws.Range("H7").Value = 3.744116
$ws.Range("I7").Value = 0.03084428108685718
$ws.Range("J7").Value = 0.03084428108685716
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04027333333333333
$ws.Range("N7").Value = 0.12082
$ws.Range("O7").Value = 0.8765171465674219
$ws.Range("P7").Value = 0.876517146567422
$ws.Range("Q7").Value = 0.05026267723555556
$ws.Range("R7").Value = 0.45236409512
$ws.Range("S7").Value = 0.02703554124617555
$ws.Range("T7").Value = 0.02703554124617554
$ws.Rows("8:10").Delete()
